$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Rows 2-43: refreshed Price (D) / Volume(1h) (E) quotes for unchanged coin rows
$ws.Cells.Item(2, 4).Value = '19.880.09'
$ws.Cells.Item(2, 5).Value = '  -8.57%  '
$ws.Cells.Item(3, 4).Value = '1.398.03'
$ws.Cells.Item(3, 5).Value = '  -9.29%  '
Set-TextCell 4 4 '1.004'
$ws.Cells.Item(4, 5).Value = '  +0.55%  '
Set-TextCell 5 4 '1.004'
$ws.Cells.Item(5, 5).Value = '  +0.49%  '
Set-TextCell 6 4 '273.13'
$ws.Cells.Item(6, 5).Value = '  -5.93%  '
Set-TextCell 7 4 '0.3670'
$ws.Cells.Item(7, 5).Value = '  -7.05%  '
Set-TextCell 8 4 '0.3105'
$ws.Cells.Item(8, 5).Value = '  -3.20%  '
Set-TextCell 9 4 '39.55'
$ws.Cells.Item(9, 5).Value = '  -8.36%  '
Set-TextCell 10 4 '1.002'
$ws.Cells.Item(10, 5).Value = '  -7.87%  '
Set-TextCell 11 4 '0.06467'
$ws.Cells.Item(11, 5).Value = '  -10.39%  '
Set-TextCell 12 4 '1.004'
$ws.Cells.Item(12, 5).Value = '  +0.57%  '
Set-TextCell 13 4 '5.414'
$ws.Cells.Item(13, 5).Value = '  -6.36%  '
Set-TextCell 14 4 '17.22'
$ws.Cells.Item(14, 5).Value = '  -7.01%  '
Set-TextCell 15 4 '6.118'
$ws.Cells.Item(15, 5).Value = '  -8.18%  '
$ws.Cells.Item(16, 4).Value = '1.397.63'
$ws.Cells.Item(16, 5).Value = '  -9.65%  '
Set-TextCell 17 4 '0.00001009'
$ws.Cells.Item(17, 5).Value = '  -8.59%  '
Set-TextCell 18 4 '0.05680'
$ws.Cells.Item(18, 5).Value = '  -14.12%  '
Set-TextCell 19 4 '1.004'
$ws.Cells.Item(19, 5).Value = '  +0.56%  '
Set-TextCell 20 4 '69.92'
$ws.Cells.Item(20, 5).Value = '  -17.18%  '
$ws.Cells.Item(21, 5).Value = '  -10.20%  '
Set-TextCell 22 4 '14.60'
$ws.Cells.Item(22, 5).Value = '  -6.75%  '
Set-TextCell 23 4 '10.98'
$ws.Cells.Item(23, 5).Value = '  +0.53%  '
Set-TextCell 24 4 '2.267'
$ws.Cells.Item(24, 5).Value = '  -4.35%  '
$ws.Cells.Item(25, 4).Value = '19.875.82'
$ws.Cells.Item(25, 5).Value = '  -8.59%  '
Set-TextCell 26 4 '2.206'
$ws.Cells.Item(26, 5).Value = '  -8.82%  '
Set-TextCell 27 4 '134.97'
$ws.Cells.Item(27, 5).Value = '  -11.39%  '
Set-TextCell 28 4 '16.83'
$ws.Cells.Item(28, 5).Value = '  -9.42%  '
$ws.Cells.Item(29, 4).Value = '1.557.44'
$ws.Cells.Item(29, 5).Value = '  -9.16%  '
Set-TextCell 30 4 '108.91'
$ws.Cells.Item(30, 5).Value = '  -7.62%  '
Set-TextCell 31 4 '4.070'
$ws.Cells.Item(31, 5).Value = '  -16.53%  '
Set-TextCell 32 4 '5.249'
$ws.Cells.Item(32, 5).Value = '  -14.89%  '
Set-TextCell 33 4 '0.8051'
$ws.Cells.Item(33, 5).Value = '  -17.96%  '
Set-TextCell 34 4 '0.07626'
$ws.Cells.Item(34, 5).Value = '  -6.52%  '
Set-TextCell 35 4 '8.339'
$ws.Cells.Item(35, 5).Value = '  -3.60%  '
Set-TextCell 36 4 '1.426'
$ws.Cells.Item(36, 5).Value = '  -4.48%  '
Set-TextCell 37 4 '0.05754'
$ws.Cells.Item(37, 5).Value = '  -4.85%  '
Set-TextCell 38 4 '4.787'
$ws.Cells.Item(38, 5).Value = '  -8.72%  '
Set-TextCell 39 4 '1.003'
$ws.Cells.Item(39, 5).Value = '  +0.47%  '
Set-TextCell 40 4 '0.02054'
$ws.Cells.Item(40, 5).Value = '  -8.96%  '
Set-TextCell 41 4 '0.1886'
$ws.Cells.Item(41, 5).Value = '  -8.33%  '
Set-TextCell 42 4 '10.27'
$ws.Cells.Item(42, 5).Value = '  -10.12%  '
Set-TextCell 43 4 '1.091'
$ws.Cells.Item(43, 5).Value = '  -8.19%  '

# Rows 44-46: ranking reorder (PancakeSwap moved up; TheSandbox/EnergySwap shifted down)
$ws.Cells.Item(44, 2).Value = 'PancakeSwap'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 44 4 '3.507'
$ws.Cells.Item(44, 5).Value = '  -6.18%  '
$ws.Cells.Item(45, 2).Value = 'TheSandbox'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell 45 4 '0.5233'
$ws.Cells.Item(45, 5).Value = '  -10.81%  '
$ws.Cells.Item(46, 2).Value = 'EnergySwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 46 4 '12.17'
$ws.Cells.Item(46, 5).Value = '  -8.00%  '

# Rows 47-51: refreshed Price (D) / Volume(1h) (E) quotes for unchanged coin rows
Set-TextCell 47 4 '0.5063'
$ws.Cells.Item(47, 5).Value = '  -9.89%  '
Set-TextCell 48 4 '111.09'
$ws.Cells.Item(48, 5).Value = '  -5.29%  '
Set-TextCell 49 4 '1.747'
$ws.Cells.Item(49, 5).Value = '  -8.56%  '
Set-TextCell 50 4 '1.007'
$ws.Cells.Item(50, 5).Value = '  +0.74%  '
Set-TextCell 51 4 '1.032'
$ws.Cells.Item(51, 5).Value = '  -12.08%  '
